$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 21 (2025Q3) metrics per updated source data
$ws.Range("C21").Value = 171
$ws.Range("D21").Value = 155
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 44.41260744985674
